$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4667.114
$ws.Range("I15").Value = 4667.114
$ws.Range("K15").Value = 14001.342
$ws.Range("M15").Value = -13832.342
$ws.Range("H107").Value = 6443.6113
$ws.Range("J107").Value = 1343
$ws.Range("L107").Value = 1343
$ws.Range("N107").Value = -5183
$ws.Range("H111").Value = 3116
$ws.Range("I111").Value = 3000
$ws.Range("J111").Value = 3232
$ws.Range("K111").Value = 9000
$ws.Range("L111").Value = 9696
$ws.Range("M111").Value = -5933
$ws.Range("N111").Value = -15830
$ws.Range("H121").Value = 1500.0454
$ws.Range("J121").Value = 1500.0454
$ws.Range("L121").Value = 4500.1362
$ws.Range("N121").Value = -7994.1362
$ws.Range("H125").Value = 2013
$ws.Range("I125").Value = 1659.2
$ws.Range("J125").Value = 2366.8
$ws.Range("K125").Value = 14932.8
$ws.Range("L125").Value = 21301.2
$ws.Range("M125").Value = -12472.8
$ws.Range("N125").Value = -26221.2
$ws.Range("H132").Value = 15831.585
$ws.Range("I132").Value = 2235.8333
$ws.Range("K132").Value = 6707.499899999999
$ws.Range("M132").Value = -4177.499899999999
$ws.Range("H137").Value = 3216490.2
$ws.Range("I137").Value = 7006397
$ws.Range("J137").Value = 9646
$ws.Range("K137").Value = 21019191
$ws.Range("L137").Value = 28938
$ws.Range("M137").Value = -21016641
$ws.Range("N137").Value = -34038
$ws.Range("H138").Value = 2277.5317
$ws.Range("I138").Value = 1422.4839
$ws.Range("J138").Value = 2829.75
$ws.Range("K138").Value = 4267.4517
$ws.Range("L138").Value = 8489.25
$ws.Range("M138").Value = 872.5483000000004
$ws.Range("N138").Value = -18769.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1763.4828
$ws.Range("I2").Value = 1957.238
$ws.Range("J2").Value = 1254.875
$ws.Range("K2").Value = 1957.238
$ws.Range("L2").Value = 1254.875
$ws.Range("M2").Value = -1844.238
$ws.Range("N2").Value = -1480.875
$ws.Range("H32").Value = 12313.695
$ws.Range("I32").Value = 11141.132
$ws.Range("J32").Value = 22671.334
$ws.Range("K32").Value = 11141.132
$ws.Range("L32").Value = 22671.334
$ws.Range("M32").Value = -10854.132
$ws.Range("N32").Value = -23245.334
$ws.Range("H46").Value = 8602.625
$ws.Range("I46").Value = 8650
$ws.Range("J46").Value = 8586.833000000001
$ws.Range("K46").Value = 8650
$ws.Range("L46").Value = 8586.833000000001
$ws.Range("M46").Value = -8331
$ws.Range("N46").Value = -9224.833000000001
$ws.Range("H61").Value = 1493.7693
$ws.Range("I61").Value = 982.075
$ws.Range("J61").Value = 3199.4167
$ws.Range("K61").Value = 982.075
$ws.Range("L61").Value = 3199.4167
$ws.Range("M61").Value = -770.075
$ws.Range("N61").Value = -3623.4167
$ws.Range("H110").Value = 1695.55
$ws.Range("I110").Value = 1556.1666
$ws.Range("J110").Value = 2950
$ws.Range("K110").Value = 1556.1666
$ws.Range("L110").Value = 2950
$ws.Range("M110").Value = 488.8334
$ws.Range("N110").Value = -7040
$ws.Range("H116").Value = 1763.4828
$ws.Range("I116").Value = 1957.238
$ws.Range("J116").Value = 1254.875
$ws.Range("K116").Value = 1957.238
$ws.Range("L116").Value = 1254.875
$ws.Range("M116").Value = 336.7619999999999
$ws.Range("N116").Value = -5842.875
$ws.Range("H136").Value = 1493.7693
$ws.Range("I136").Value = 982.075
$ws.Range("J136").Value = 3199.4167
$ws.Range("K136").Value = 2946.225
$ws.Range("L136").Value = 9598.250100000001
$ws.Range("M136").Value = -396.2250000000004
$ws.Range("N136").Value = -14698.2501
$ws.Range("H138").Value = 57942.668
$ws.Range("J138").Value = 57942.668
$ws.Range("L138").Value = 57942.668
$ws.Range("N138").Value = -68222.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1763.4828
$ws.Range("I3").Value = 1957.238
$ws.Range("J3").Value = 1254.875
$ws.Range("K3").Value = 1957.238
$ws.Range("L3").Value = 1254.875
$ws.Range("M3").Value = -1843.238
$ws.Range("N3").Value = -1482.875
$ws.Range("H134").Value = 2035.8918
$ws.Range("I134").Value = 1600.9688
$ws.Range("J134").Value = 4819.4
$ws.Range("K134").Value = 4802.9064
$ws.Range("L134").Value = 14458.2
$ws.Range("M134").Value = -2267.9064
$ws.Range("N134").Value = -19528.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3666521
$ws.Range("I31").Value = 1592.8334
$ws.Range("J31").Value = 6807888
$ws.Range("K31").Value = 1592.8334
$ws.Range("L31").Value = 6807888
$ws.Range("M31").Value = -1297.8334
$ws.Range("N31").Value = -6808478
$ws.Range("H34").Value = 3666521
$ws.Range("I34").Value = 1592.8334
$ws.Range("J34").Value = 6807888
$ws.Range("K34").Value = 1592.8334
$ws.Range("L34").Value = 6807888
$ws.Range("M34").Value = -1390.8334
$ws.Range("N34").Value = -6808292
$ws.Range("H58").Value = 2453.2856
$ws.Range("I58").Value = 1545.3846
$ws.Range("J58").Value = 3240.1333
$ws.Range("K58").Value = 1545.3846
$ws.Range("L58").Value = 3240.1333
$ws.Range("M58").Value = -1342.3846
$ws.Range("N58").Value = -3646.1333
$ws.Range("H99").Value = 1984.7693
$ws.Range("J99").Value = 2037.8182
$ws.Range("L99").Value = 2037.8182
$ws.Range("N99").Value = -5033.8182
$ws.Range("H107").Value = 599.9259
$ws.Range("I107").Value = 408
$ws.Range("J107").Value = 1271.6666
$ws.Range("K107").Value = 408
$ws.Range("L107").Value = 1271.6666
$ws.Range("M107").Value = 1512
$ws.Range("N107").Value = -5111.6666
$ws.Range("H112").Value = 51998
$ws.Range("J112").Value = 51998
$ws.Range("L112").Value = 51998
$ws.Range("N112").Value = -54952
$ws.Range("H126").Value = 1984.7693
$ws.Range("J126").Value = 2037.8182
$ws.Range("L126").Value = 6113.4546
$ws.Range("N126").Value = -11053.4546
$ws.Range("H134").Value = 760245.8
$ws.Range("I134").Value = 455357.47
$ws.Range("J134").Value = 2335502
$ws.Range("K134").Value = 1366072.41
$ws.Range("L134").Value = 7006506
$ws.Range("M134").Value = -1363537.41
$ws.Range("N134").Value = -7011576
$ws.Range("H136").Value = 2453.2856
$ws.Range("I136").Value = 1545.3846
$ws.Range("J136").Value = 3240.1333
$ws.Range("K136").Value = 4636.1538
$ws.Range("L136").Value = 9720.3999
$ws.Range("M136").Value = -2086.1538
$ws.Range("N136").Value = -14820.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 37931.25
$ws.Range("I18").Value = 43307.145
$ws.Range("K18").Value = 129921.435
$ws.Range("M18").Value = -129752.435
$ws.Range("H34").Value = 1684.7693
$ws.Range("J34").Value = 2055.2
$ws.Range("L34").Value = 6165.599999999999
$ws.Range("N34").Value = -6333.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3925.2
$ws.Range("I107").Value = 465
$ws.Range("K107").Value = 465
$ws.Range("M107").Value = 1455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4267.923
$ws.Range("I40").Value = 3400.1428
$ws.Range("J40").Value = 5280.3335
$ws.Range("K40").Value = 3400.1428
$ws.Range("L40").Value = 5280.3335
$ws.Range("M40").Value = -3264.1428
$ws.Range("N40").Value = -5552.3335
$ws.Range("H61").Value = 2106.375
$ws.Range("I61").Value = 1754.7273
$ws.Range("J61").Value = 2880
$ws.Range("K61").Value = 1754.7273
$ws.Range("L61").Value = 2880
$ws.Range("M61").Value = -1552.7273
$ws.Range("N61").Value = -3284
$ws.Range("H113").Value = 2106.375
$ws.Range("I113").Value = 1754.7273
$ws.Range("J113").Value = 2880
$ws.Range("K113").Value = 1754.7273
$ws.Range("L113").Value = 2880
$ws.Range("M113").Value = 415.2727
$ws.Range("N113").Value = -7220
$ws.Range("H122").Value = 69040.266
$ws.Range("I122").Value = 93127.63
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 279382.89
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -276932.89
$ws.Range("N122").Value = -13300
$ws.Range("H132").Value = 3463.4146
$ws.Range("I132").Value = 2820.8635
$ws.Range("K132").Value = 8462.5905
$ws.Range("M132").Value = -5932.5905
$ws.Range("H135").Value = 32446.924
$ws.Range("J135").Value = 32446.924
$ws.Range("L135").Value = 32446.924
$ws.Range("N135").Value = -42586.924

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2826.6
$ws.Range("I62").Value = 2777.6667
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 2777.6667
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -2153.6667
$ws.Range("N62").Value = -4148
$ws.Range("H65").Value = 2826.6
$ws.Range("I65").Value = 2777.6667
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 13888.3335
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -10768.3335
$ws.Range("N65").Value = -20740
$ws.Range("H107").Value = 4762879.5
$ws.Range("I107").Value = 910.8
$ws.Range("J107").Value = 16667800
$ws.Range("K107").Value = 2732.4
$ws.Range("L107").Value = 50003400
$ws.Range("M107").Value = -812.3999999999996
$ws.Range("N107").Value = -50007240
